$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new commission period ("2509") was billed in addition to the existing
# 2507 / 2508 rows, so duplicate the last data row (row 17, period 2508)
# down into a new row 18, keeping the same worker/account details and
# amounts, then relabel the period.
$ws.Rows(18).Insert()
$ws.Range("B17:J17").Copy($ws.Range("B18:J18"))

# Row 17 is no longer the last row of the table, so it should take on the
# "interior" row formatting that row 16 already uses (the new row 18 keeps
# the old "last row" formatting, e.g. the outer bottom border).
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Label the newly added row with its own period.
$ws.Range("E18").Value = "2509"

# Refresh the summary figures: one more period of mora for the same
# worker, and one more period counted overall.
$ws.Range("E11").Value = 170820
$ws.Range("F13").Value = 3
